# Edit script: update "DataSource - Emision Motor" worksheet
# Commit: "se arregl\u00f3 el error con el item PantallaInicio, esto afectpo el
#  recording Personas.rxrec Se modif los Datasources para hacer la regre
#  en PreProd R30"
#
# Refreshes the existing data row (NroCuenta / FechaInicio / SumaAsegurada /
# Patente / Motor / Chasis) and adds two more sample rows corresponding to
# other MetodoDePago values (Tarjeta de Credito, Debito Bancario).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-CellFormat($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Row 2: refresh existing record values (formats/styles stay as they are)
# ---------------------------------------------------------------------
$ws.Range("F2").Value = 2617100594
$ws.Range("J2").Value = "'10/03/2021"
$ws.Range("W2").Value = 1380000
$ws.Range("Y2").Value = "RPM001"
$ws.Range("Z2").Value = "ABC12SRPM001"
$ws.Range("AA2").Value = "ZAZ123SRPM001"

# ---------------------------------------------------------------------
# Row 3: new record (Tarjeta de Credito)
# ---------------------------------------------------------------------
# columns that need a non-default cell style get format copied from row 2
# (single cell at a time, so unrelated blank cells are left untouched)
Copy-CellFormat "B2" "B3"
Copy-CellFormat "C2" "C3"
Copy-CellFormat "F2" "F3"
Copy-CellFormat "J2" "J3"
Copy-CellFormat "B2" "O3"
Copy-CellFormat "B2" "P3"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "'preproducciongestion.segurossura.com.ar"
$ws.Range("C3").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("D3").Value = "su"
$ws.Range("E3").Value = "silverarrow"
$ws.Range("F3").Value = 2617100594
$ws.Range("G3").Value = "Motor"
$ws.Range("H3").Value = "Menos de 5 vehículos"
$ws.Range("I3").Value = "Anual"
$ws.Range("J3").Value = "'10/03/2021"
$ws.Range("K3").Value = "Tarjeta de Crédito"
$ws.Range("L3").Value = "No"
$ws.Range("N3").Value = "Naranja"
$ws.Range("O3").Value = "'5895627074474000"
$ws.Range("P3").Value = "'03/2025"
$ws.Range("S3").Value = 2021
$ws.Range("T3").Value = "si"
$ws.Range("U3").Value = "CHEVROLET"
$ws.Range("V3").Value = "ONIX 1.2 L/19"
$ws.Range("W3").Value = 1380000
$ws.Range("X3").Value = "CPremium - Resp. Civil-Robo/Incendio Total y Parcial Daños Totales por Accidente"
$ws.Range("Y3").Value = "RPM002"
$ws.Range("Z3").Value = "ABC12SRPM002"
$ws.Range("AA3").Value = "ZAZ123SRPM002"
$ws.Range("AB3").Value = "Sin Actividad"
$ws.Range("AC3").Value = "No"

# ---------------------------------------------------------------------
# Row 4: new record (Debito Bancario)
# ---------------------------------------------------------------------
Copy-CellFormat "B2" "B4"
Copy-CellFormat "C2" "C4"
Copy-CellFormat "F2" "F4"
Copy-CellFormat "J2" "J4"
Copy-CellFormat "B2" "R4"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "'preproducciongestion.segurossura.com.ar"
$ws.Range("C4").Value = "https://preproducciongestion.segurossura.com.ar/pc/PolicyCenter.do"
$ws.Range("D4").Value = "su"
$ws.Range("E4").Value = "silverarrow"
$ws.Range("F4").Value = 2617100594
$ws.Range("G4").Value = "Motor"
$ws.Range("H4").Value = "Menos de 5 vehículos"
$ws.Range("I4").Value = "Anual"
$ws.Range("J4").Value = "'10/03/2021"
$ws.Range("K4").Value = "Débito Bancario"
$ws.Range("L4").Value = "No"
$ws.Range("Q4").Value = "SNP"
$ws.Range("R4").Value = "'0340010400044915199004"
$ws.Range("S4").Value = 2021
$ws.Range("T4").Value = "si"
$ws.Range("U4").Value = "CHEVROLET"
$ws.Range("V4").Value = "ONIX 1.2 L/19"
$ws.Range("W4").Value = 1380000
$ws.Range("X4").Value = "CPremium - Resp. Civil-Robo/Incendio Total y Parcial Daños Totales por Accidente"
$ws.Range("Y4").Value = "RPM003"
$ws.Range("Z4").Value = "ABC12SRPM003"
$ws.Range("AA4").Value = "ZAZ123SRPM003"
$ws.Range("AB4").Value = "Sin Actividad"
$ws.Range("AC4").Value = "No"

# ---------------------------------------------------------------------
# Sheet view: scroll right so column Q is the left-most visible column,
# with Y2 as the active selection.
# ---------------------------------------------------------------------
$ws.Range("Y2").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("Q1").Column
